{"js": "// Update the date line and each three-digit x one-digit multiplication\n// expression in the practice table to the new values from the latest\n// generated worksheet.\nconst replacements = [\n  [\"2024-08-20 Tuesday\", \"2024-08-21 Wednesday\"],\n  [\"558\u00d73=1674\", \"492\u00d73=1476\"],\n  [\"788\u00d77=5516\", \"889\u00d76=5334\"],\n  [\"733\u00d75=3665\", \"250\u00d73=750\"],\n  [\"152\u00d74=608\", \"503\u00d73=1509\"],\n  [\"377\u00d76=2262\", \"173\u00d78=1384\"],\n  [\"973\u00d76=5838\", \"314\u00d78=2512\"],\n  [\"201\u00d77=1407\", \"685\u00d77=4795\"],\n  [\"483\u00d79=4347\", \"122\u00d79=1098\"],\n  [\"512\u00d79=4608\", \"830\u00d72=1660\"],\n  [\"479\u00d77=3353\", \"158\u00d77=1106\"],\n  [\"101\u00d73=303\", \"508\u00d78=4064\"],\n  [\"157\u00d75=785\", \"191\u00d79=1719\"],\n  [\"671\u00d76=4026\", \"869\u00d79=7821\"],\n  [\"612\u00d74=2448\", \"657\u00d77=4599\"],\n  [\"799\u00d73=2397\", \"377\u00d74=1508\"],\n  [\"231\u00d78=1848\", \"880\u00d74=3520\"],\n  [\"540\u00d72=1080\", \"928\u00d73=2784\"],\n  [\"892\u00d77=6244\", \"777\u00d72=1554\"],\n  [\"720\u00d72=1440\", \"790\u00d76=4740\"],\n  [\"867\u00d73=2601\", \"719\u00d79=6471\"],\n  [\"177\u00d78=1416\", \"572\u00d78=4576\"],\n  [\"233\u00d75=1165\", \"621\u00d79=5589\"],\n  [\"293\u00d72=586\", \"915\u00d72=1830\"],\n  [\"903\u00d73=2709\", \"421\u00d77=2947\"],\n  [\"766\u00d76=4596\", \"451\u00d72=902\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and each three-digit x one-digit multiplication\n# expression in the practice table to the new values from the latest\n# generated worksheet.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Old=\"2024-08-20 Tuesday\"; New=\"2024-08-21 Wednesday\"},\n    @{Old=\"558\u00d73=1674\"; New=\"492\u00d73=1476\"},\n    @{Old=\"788\u00d77=5516\"; New=\"889\u00d76=5334\"},\n    @{Old=\"733\u00d75=3665\"; New=\"250\u00d73=750\"},\n    @{Old=\"152\u00d74=608\"; New=\"503\u00d73=1509\"},\n    @{Old=\"377\u00d76=2262\"; New=\"173\u00d78=1384\"},\n    @{Old=\"973\u00d76=5838\"; New=\"314\u00d78=2512\"},\n    @{Old=\"201\u00d77=1407\"; New=\"685\u00d77=4795\"},\n    @{Old=\"483\u00d79=4347\"; New=\"122\u00d79=1098\"},\n    @{Old=\"512\u00d79=4608\"; New=\"830\u00d72=1660\"},\n    @{Old=\"479\u00d77=3353\"; New=\"158\u00d77=1106\"},\n    @{Old=\"101\u00d73=303\"; New=\"508\u00d78=4064\"},\n    @{Old=\"157\u00d75=785\"; New=\"191\u00d79=1719\"},\n    @{Old=\"671\u00d76=4026\"; New=\"869\u00d79=7821\"},\n    @{Old=\"612\u00d74=2448\"; New=\"657\u00d77=4599\"},\n    @{Old=\"799\u00d73=2397\"; New=\"377\u00d74=1508\"},\n    @{Old=\"231\u00d78=1848\"; New=\"880\u00d74=3520\"},\n    @{Old=\"540\u00d72=1080\"; New=\"928\u00d73=2784\"},\n    @{Old=\"892\u00d77=6244\"; New=\"777\u00d72=1554\"},\n    @{Old=\"720\u00d72=1440\"; New=\"790\u00d76=4740\"},\n    @{Old=\"867\u00d73=2601\"; New=\"719\u00d79=6471\"},\n    @{Old=\"177\u00d78=1416\"; New=\"572\u00d78=4576\"},\n    @{Old=\"233\u00d75=1165\"; New=\"621\u00d79=5589\"},\n    @{Old=\"293\u00d72=586\"; New=\"915\u00d72=1830\"},\n    @{Old=\"903\u00d73=2709\"; New=\"421\u00d77=2947\"},\n    @{Old=\"766\u00d76=4596\"; New=\"451\u00d72=902\"}\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)\n}\n"}
